# Update LR-pair TPM-derived statistics in Sheet1 to reflect the new TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.104631
$ws.Range("H2").Value = 3.313893
$ws.Range("I2").Value = 0.8734451962653081
$ws.Range("J2").Value = 0.8734451962653083
$ws.Range("M2").Value = 19.98610666666667
$ws.Range("N2").Value = 59.95832
$ws.Range("O2").Value = 0.2969043109767812
$ws.Range("P2").Value = 0.2969043109767812
$ws.Range("Q2").Value = 22.07727299330667
$ws.Range("R2").Value = 198.69545693976
$ws.Range("S2").Value = 0.2593296441731308
$ws.Range("T2").Value = 0.2593296441731308

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.104631
$ws.Range("H3").Value = 3.313893
$ws.Range("I3").Value = 0.8734451962653081
$ws.Range("J3").Value = 0.8734451962653083
$ws.Range("O3").Value = 0.4664722083712238
$ws.Range("P3").Value = 0.4664722083712239
$ws.Range("Q3").Value = 34.68603825293567
$ws.Range("R3").Value = 312.174344276421
$ws.Range("S3").Value = 0.4074379095931153
$ws.Range("T3").Value = 0.4074379095931154

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.104631
$ws.Range("H4").Value = 3.313893
$ws.Range("I4").Value = 0.8734451962653081
$ws.Range("J4").Value = 0.8734451962653083
$ws.Range("O4").Value = 0.236623480651995
$ws.Range("P4").Value = 0.236623480651995
$ws.Range("Q4").Value = 17.594898375824
$ws.Range("R4").Value = 158.354085382416
$ws.Range("S4").Value = 0.2066776424990621
$ws.Range("T4").Value = 0.2066776424990621

# Row 5
$ws.Range("I5").Value = 0.1265548037346918
$ws.Range("J5").Value = 0.1265548037346918
$ws.Range("M5").Value = 19.98610666666667
$ws.Range("N5").Value = 59.95832
$ws.Range("O5").Value = 0.2969043109767812
$ws.Range("P5").Value = 0.2969043109767812
$ws.Range("Q5").Value = 3.198809682177778
$ws.Range("R5").Value = 28.7892871396
$ws.Range("S5").Value = 0.03757466680365044
$ws.Range("T5").Value = 0.03757466680365044

# Row 6
$ws.Range("I6").Value = 0.1265548037346918
$ws.Range("J6").Value = 0.1265548037346918
$ws.Range("O6").Value = 0.4664722083712238
$ws.Range("P6").Value = 0.4664722083712239
$ws.Range("S6").Value = 0.05903429877810848
$ws.Range("T6").Value = 0.05903429877810849

# Row 7
$ws.Range("I7").Value = 0.1265548037346918
$ws.Range("J7").Value = 0.1265548037346918
$ws.Range("O7").Value = 0.236623480651995
$ws.Range("P7").Value = 0.236623480651995
$ws.Range("S7").Value = 0.02994583815293286
$ws.Range("T7").Value = 0.02994583815293287
